$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new blank row above the old row 4 (becomes row 3). Excel shifts
#    every row from 4..16 down to 5..17, carrying cell values/styles/heights
#    with them automatically.
# ---------------------------------------------------------------------------
$ws.Rows("3").Insert()

# ---------------------------------------------------------------------------
# 2) Stash the formatting of the five cells that carry hyperlinks, because
#    re-creating the hyperlinks (needed since the engine does not re-point
#    existing Hyperlinks ranges on a row insert) resets cell formatting to
#    the default "Hyperlink" look. We copy format-only into scratch cells far
#    off to the side and paste them back once the links are rebuilt.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("Z5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Drop every hyperlink and rebuild them pointing at the shifted cells.
#    For the two links whose cached "display" text differs from the real
#    cell text, set it via TextToDisplay (which stamps the cell value too)
#    and then restore the true cell text immediately after - the cached
#    display string on the hyperlink itself is unaffected by that follow-up
#    write. The other three links are added without TextToDisplay, which
#    leaves the cell's existing value (and, for B2, its phonetic guide)
#    completely untouched.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B5"), "https://knot-designs.com/", "", "", "https://knot-designs.com/")
$ws.Range("B5").Value2 = "Maker's Watch Knot (knot-designs.com)"

$ws.Hyperlinks.Add($ws.Range("B9"), "https://undone.co.jp/", "", "", "https://undone.co.jp/")
$ws.Range("B9").Value2 = "【UNDONE（アンダーン）】カスタマイズリストウェアブランド"

$ws.Hyperlinks.Add($ws.Range("B13"), "https://renautus.com/")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://custom-watch-mania.com/")
$ws.Hyperlinks.Add($ws.Range("B17"), "index.html")

# ---------------------------------------------------------------------------
# 4) Restore the original formatting onto the (now re-hyperlinked) cells.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("Z4").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("Z1:Z5").Clear()

# ---------------------------------------------------------------------------
# 5) Fill in the brand-new banner row (row 3): a bold red warning telling
#    internal reviewers the page isn't public yet and its links are
#    unapproved.
# ---------------------------------------------------------------------------
$banner = $ws.Range("B3")
$banner.Value2 = "本ページは社内検討用の公開前ページです。`n以下のリンクにはまだ許可をいただいておりません"

# Base the new style on the existing title cell (B2) style, then tweak it:
# bold, 16pt, red, centered + wrapped.
$ws.Range("B2").Copy()
$banner.PasteSpecial(-4122)
$banner.Font.Bold = $true
$banner.Font.Size = 16
$banner.Font.Color = 255
$banner.HorizontalAlignment = -4108
$banner.VerticalAlignment = -4108
$banner.WrapText = $true
$ws.Rows("3").RowHeight = 51

$excel.CutCopyMode = $false
